$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.224.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5137"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3903"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.119"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.241"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001108"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06655"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.249.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.091.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.509"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.044"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.897"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.31%  "

$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.773"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02462"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06559"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2197"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.212"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6528"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.035"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.233"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6123"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.292"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.681"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.022"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.38%  "
